$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.627703740713571
    "D2" = 0.1279418244843811
    "E2" = 1.049007991950532
    "F2" = 2.851904209353165
    "G2" = 0.002436819216202656
    "I2" = 0.8100945568519404
    "L2" = 0.9828555794891258
    "B3" = 1.502000301002283
    "D3" = 0.1210811433409731
    "E3" = 0.9137212164588675
    "F3" = 2.677280718396645
    "G3" = 0.002448753308349282
    "I3" = 0.8388043192501051
    "L3" = 0.8865792465441871
    "B4" = 1.425224270507726
    "D4" = 0.1170715228459045
    "E4" = 0.830645632023419
    "F4" = 2.573502035141217
    "G4" = 0.002456434227656734
    "I4" = 0.8574469461567213
    "L4" = 0.8279107672790644
    "B5" = 1.394038266786822
    "D5" = 0.1154869157118981
    "E5" = 0.7967821351280406
    "F5" = 2.532049272878936
    "G5" = 0.002459653617842334
    "I5" = 0.8652976201404297
    "L5" = 0.804110248128552
    "B6" = 1.388865909464585
    "D6" = 0.1152267251940629
    "E6" = 0.7911583445752797
    "F6" = 2.525215892337172
    "G6" = 0.002460193606535347
    "I6" = 0.8666164950110113
    "L6" = 0.8001645406273781
    "B7" = 1.424803278478805
    "D7" = 0.1170499546568635
    "E7" = 0.8301889850440887
    "F7" = 2.572939631767554
    "G7" = 0.002456477283086463
    "I7" = 0.85755179802565
    "L7" = 0.8275893565240438
    "B8" = 1.584275962571951
    "D8" = 0.1255331414052421
    "E8" = 1.002358481034975
    "F8" = 2.790963191602515
    "G8" = 0.002440861061126398
    "I8" = 0.8197821673752967
    "L8" = 0.949564126541901
    "B9" = 1.90030962925249
    "D9" = 0.1438557688730526
    "E9" = 1.340266070135186
    "F9" = 3.247089622837478
    "G9" = 0.002413018072622909
    "I9" = 0.7538349340809614
    "L9" = 1.192521911556071
    "B10" = 2.134666396050761
    "D10" = 0.1584567956451792
    "E10" = 1.589231442490103
    "F10" = 3.601478933168721
    "G10" = 0.002394223940575129
    "I10" = 0.7104261972580037
    "L10" = 1.37366559431797
    "B11" = 2.241788925065975
    "D11" = 0.1653715399536964
    "E11" = 1.702768003713715
    "F11" = 3.767302172884797
    "G11" = 0.002386027703980198
    "I11" = 0.6917941869336053
    "L11" = 1.456726261016058
    "B12" = 2.282429864075141
    "D12" = 0.1680312800319257
    "E12" = 1.745811843398258
    "F12" = 3.8307928703789
    "G12" = 0.002382974229842121
    "I12" = 0.6849008769120006
    "L12" = 1.488280196668484
    "B13" = 2.273673694098648
    "D13" = 0.167456588000789
    "E13" = 1.736539200627902
    "F13" = 3.817087455740705
    "G13" = 0.002383629623411938
    "I13" = 0.6863782283812094
    "L13" = 1.48147991483836
    "B14" = 2.245130943514368
    "D14" = 0.1655895190452554
    "E14" = 1.706308187470739
    "F14" = 3.772511421430465
    "G14" = 0.002385775488753949
    "I14" = 0.6912238077329027
    "L14" = 1.459320168879856
    "B15" = 2.227657647724754
    "D15" = 0.1644513233167686
    "E15" = 1.687797624088802
    "F15" = 3.745299121090198
    "G15" = 0.002387096422454382
    "I15" = 0.6942130495798438
    "L15" = 1.445759979979414
    "B16" = 2.127676188747103
    "D16" = 0.1580105628023034
    "E16" = 1.581818026203479
    "F16" = 3.590737648826519
    "G16" = 0.002394766647112012
    "I16" = 0.7116664389386953
    "L16" = 1.368251131946465
    "B17" = 2.066473612632137
    "D17" = 0.1541306512948779
    "E17" = 1.516881636278981
    "F17" = 3.497123811948228
    "G17" = 0.002399562188699896
    "I17" = 0.7226604119843412
    "L17" = 1.320874636389192
    "B18" = 2.031319597581728
    "D18" = 0.1519245656822648
    "E18" = 1.479557787579211
    "F18" = 3.443711391791339
    "G18" = 0.002402353745505787
    "I18" = 0.7290886710821436
    "L18" = 1.293686434685981
    "B19" = 2.019425245110881
    "D19" = 0.1511819470414082
    "E19" = 1.466924727014657
    "F19" = 3.42570000077049
    "G19" = 0.002403304652365929
    "I19" = 0.7312831139950546
    "L19" = 1.284491358169987
    "B20" = 2.072983738557582
    "D20" = 0.1545410167604189
    "E20" = 1.523791482744969
    "F20" = 3.507044232991149
    "G20" = 0.002399048253912876
    "I20" = 0.721479222963314
    "L20" = 1.325911536144474
    "B21" = 2.25351256062396
    "D21" = 0.1661367854186437
    "E21" = 1.715186332945791
    "F21" = 3.785585288025459
    "G21" = 0.00238514383639912
    "I21" = 0.6897961252717684
    "L21" = 1.465826241411719
    "B22" = 2.371942535799064
    "D22" = 0.1739568890251917
    "E22" = 1.840570189604733
    "F22" = 3.971707503006769
    "G22" = 0.002376349206692565
    "I22" = 0.6700359044085449
    "L22" = 1.557858360646094
    "B23" = 2.308692808068884
    "D23" = 0.1697603456803449
    "E23" = 1.77362006977188
    "F23" = 3.871985710926481
    "G23" = 0.002381016464591697
    "I23" = 0.6804950534900804
    "L23" = 1.508683033661214
    "B24" = 2.070040410621118
    "D24" = 0.1543554142920556
    "E24" = 1.520667514011393
    "F24" = 3.502557944180836
    "G24" = 0.002399280496231054
    "I24" = 0.7220129032725895
    "L24" = 1.323634201050254
    "B25" = 1.814445085922273
    "D25" = 0.1387071309127492
    "E25" = 1.248769112256042
    "F25" = 3.120447710381882
    "G25" = 0.002420256053532826
    "I25" = 0.7707967524908899
    "L25" = 1.126357193063086
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

